$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.274.99"
$ws.Range("E2").Value = "  -5.11%  "

$ws.Range("D3").Value = "1.673.01"
$ws.Range("E3").Value = "  -2.69%  "

$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").Value = "'217.99"
$ws.Range("E5").Value = "  -2.68%  "

$ws.Range("D6").Value = "'0.5116"
$ws.Range("E6").Value = "  -9.78%  "

$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").Value = "'0.2664"
$ws.Range("E8").Value = "  -1.04%  "

$ws.Range("D9").Value = "'0.06392"
$ws.Range("E9").Value = "  -2.24%  "

$ws.Range("D10").Value = "'21.48"
$ws.Range("E10").Value = "  -5.44%  "

$ws.Range("D11").Value = "'0.07379"
$ws.Range("E11").Value = "  -1.73%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "'4.561"
$ws.Range("E12").Value = "  -1.61%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.668.88"
$ws.Range("E13").Value = "  -3.18%  "

$ws.Range("D14").Value = "'0.5832"
$ws.Range("E14").Value = "  -1.64%  "

$ws.Range("D15").Value = "1.898.62"
$ws.Range("E15").Value = "  -2.95%  "

$ws.Range("D16").Value = "'0.000008586"
$ws.Range("E16").Value = "  +0.54%  "

$ws.Range("D17").Value = "'64.81"
$ws.Range("E17").Value = "  -11.94%  "

$ws.Range("D18").Value = "26.338.31"
$ws.Range("E18").Value = "  -4.85%  "

$ws.Range("D19").Value = "'4.957"
$ws.Range("E19").Value = "  -5.50%  "

$ws.Range("E20").Value = "  +0.29%  "

$ws.Range("D21").Value = "'10.84"
$ws.Range("E21").Value = "  -2.87%  "

$ws.Range("D22").Value = "'190.12"
$ws.Range("E22").Value = "  -4.96%  "

$ws.Range("D23").Value = "'6.220"
$ws.Range("E23").Value = "  -4.31%  "

$ws.Range("D24").Value = "'1.007"
$ws.Range("E24").Value = "  +0.23%  "

$ws.Range("D25").Value = "'143.75"
$ws.Range("E25").Value = "  -3.55%  "

$ws.Range("D26").Value = "'7.664"
$ws.Range("E26").Value = "  -3.29%  "

$ws.Range("D27").Value = "'0.1183"
$ws.Range("E27").Value = "  -2.16%  "

$ws.Range("D28").Value = "'15.66"
$ws.Range("E28").Value = "  -1.93%  "

$ws.Range("D29").Value = "'0.05889"
$ws.Range("E29").Value = "  -3.71%  "

$ws.Range("D30").Value = "'1.276"
$ws.Range("E30").Value = "  -7.27%  "

$ws.Range("D31").Value = "'1.326"
$ws.Range("E31").Value = "  -3.85%  "

$ws.Range("D32").Value = "'3.518"
$ws.Range("E32").Value = "  -4.35%  "

$ws.Range("D33").Value = "'3.526"
$ws.Range("E33").Value = "  -4.01%  "

$ws.Range("D34").Value = "'1.642"
$ws.Range("E34").Value = "  -1.00%  "

$ws.Range("D35").Value = "'1.014"
$ws.Range("E35").Value = "  -0.85%  "

$ws.Range("D36").Value = "'0.6021"
$ws.Range("E36").Value = "  -5.68%  "

$ws.Range("D37").Value = "'2.360"
$ws.Range("E37").Value = "  -2.63%  "

$ws.Range("D38").Value = "'2.650"
$ws.Range("E38").Value = "  -0.89%  "

$ws.Range("D39").Value = "'0.01619"
$ws.Range("E39").Value = "  -1.87%  "

$ws.Range("D40").Value = "'6.043"
$ws.Range("E40").Value = "  -1.38%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'0.8684"
$ws.Range("E41").Value = "  -0.33%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.074.02"
$ws.Range("E42").Value = "  -3.33%  "

$ws.Range("E43").Value = "  +0.52%  "

$ws.Range("D44").Value = "'99.73"
$ws.Range("E44").Value = "  +0.59%  "

$ws.Range("D45").Value = "1.819.80"
$ws.Range("E45").Value = "  -2.71%  "

$ws.Range("D46").Value = "'0.00000000113"
$ws.Range("E46").Value = "  +3.13%  "

$ws.Range("D47").Value = "'55.98"
$ws.Range("E47").Value = "  -4.37%  "

$ws.Range("D48").Value = "'1.008"
$ws.Range("E48").Value = "  +0.77%  "

$ws.Range("D49").Value = "'8.051"
$ws.Range("E49").Value = "  -1.73%  "

$ws.Range("E50").Value = "  -2.43%  "

$ws.Range("D51").Value = "'0.05185"
$ws.Range("E51").Value = "  -3.30%  "

Write-Host "Applied updates"